$d = $word.ActiveDocument

# Step 1: change the single space run " " before the ellipsis run to "ที่ "
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Text = "ตาราง "
$find1.Replacement.ClearFormatting()
$find1.Replacement.Text = "ตารางที่ "
$find1.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)

# Step 2: change the ellipsis run "… " to "1 "
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "… "
$find2.Replacement.ClearFormatting()
$find2.Replacement.Text = "1 "
$find2.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)
